# Raw and Clean Data from SSA for June 20th and 21th
#
# The original sheet had an extra leading "index" column (A) that is no
# longer wanted: Fecha/Confirmados/Negativos/Sospechosos/Defunciones/
# Porcentaje hospitalizados used to live in columns B:G. We drop that
# index column (shifting B:G left into A:F) and append two new rows of
# data for 2020-06-20 and 2020-06-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the old index column (A); this shifts columns B:G to A:F,
# including headers, values and per-cell styles/number formats.
$ws.Columns("A").Delete()

# After the shift, the last existing data row (20) currently carries the
# special "final row" date style (date-only, no time) that used to live on
# the old last row. Since we are about to append two more rows, that
# special style needs to move to the new final row (22) instead, and row
# 20 should fall back to the regular date/time style used by the other
# data rows (like row 19).

# Step 2: copy the special "final row" style from A20 onto the new last
# row (A22) before it gets overwritten.
$ws.Range("A20").Copy()
$ws.Range("A22").PasteSpecial(-4122)

# Step 3: copy the regular data-row style from A19 onto A20 (no longer the
# last row) and onto the new row A21.
$ws.Range("A19").Copy()
$ws.Range("A20:A21").PasteSpecial(-4122)

# Step 4: populate the two new rows of SSA data.
# Row 21 -> 2020-06-20 (serial date 44002)
$ws.Range("A21").Value2 = 44002
$ws.Range("B21").Value2 = 175202
$ws.Range("C21").Value2 = 238129
$ws.Range("D21").Value2 = 60621
$ws.Range("E21").Value2 = 20781
$ws.Range("F21").Value2 = 31.46

# Row 22 -> 2020-06-21 (serial date 44003)
$ws.Range("A22").Value2 = 44003
$ws.Range("B22").Value2 = 180545
$ws.Range("C22").Value2 = 242393
$ws.Range("D22").Value2 = 56590
$ws.Range("E22").Value2 = 21825
$ws.Range("F22").Value2 = 31.61
